$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

Set-TextCell $ws "D2" "51.954.14"
Set-TextCell $ws "E2" "  +0.38%  "
Set-TextCell $ws "D3" "2.790.36"
Set-TextCell $ws "E3" "  -0.95%  "
Set-TextCell $ws "E4" "  -0.06%  "
Set-TextCell $ws "D5" "358.26"
Set-TextCell $ws "D6" "109.32"
Set-TextCell $ws "E6" "  -3.22%  "
Set-TextCell $ws "E7" "  -0.15%  "
Set-TextCell $ws "E8" "  -0.01%  "
Set-TextCell $ws "D9" "0.592"
Set-TextCell $ws "E9" "  -1.15%  "
Set-TextCell $ws "D10" "40.06"
Set-TextCell $ws "E10" "  -3.43%  "
Set-TextCell $ws "D11" "0.0855"
Set-TextCell $ws "E11" "  +0.30%  "
Set-TextCell $ws "E12" "  +1.20%  "
Set-TextCell $ws "D13" "19.49"
Set-TextCell $ws "E13" "  -2.02%  "
Set-TextCell $ws "D14" "7.58"
Set-TextCell $ws "E14" "  -2.14%  "
Set-TextCell $ws "D15" "3.226.37"
Set-TextCell $ws "E15" "  -1.07%  "
Set-TextCell $ws "D16" "2.783.89"
Set-TextCell $ws "E16" "  -1.26%  "
Set-TextCell $ws "D17" "0.950"
Set-TextCell $ws "E17" "  +7.31%  "
Set-TextCell $ws "D18" "51.883.55"
Set-TextCell $ws "E18" "  +0.38%  "
Set-TextCell $ws "E19" "  -1.03%  "
Set-TextCell $ws "E20" "  -1.51%  "
Set-TextCell $ws "E21" "  -2.84%  "
Set-TextCell $ws "D22" "0.0₃0981"
Set-TextCell $ws "E22" "  -1.19%  "
Set-TextCell $ws "D23" "274.32"
Set-TextCell $ws "E23" "  +1.44%  "
Set-TextCell $ws "D24" "70.27"
Set-TextCell $ws "E24" "  +0.85%  "
Set-TextCell $ws "D25" "2.73"
Set-TextCell $ws "E25" "  -0.53%  "
Set-TextCell $ws "D26" "26.69"
Set-TextCell $ws "E26" "  -0.06%  "
Set-TextCell $ws "E27" "  +0.11%  "
Set-TextCell $ws "E28" "  -1.42%  "
Set-TextCell $ws "B29" "Toncoin"
Set-TextCell $ws "C29" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D29" "2.28"
Set-TextCell $ws "E29" "  +1.67%  "
Set-TextCell $ws "B30" "Kaspa"
Set-TextCell $ws "C30" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws "D30" "0.145"
Set-TextCell $ws "E30" "  +3.93%  "
Set-TextCell $ws "D31" "0.0464"
Set-TextCell $ws "E31" "  +2.77%  "
Set-TextCell $ws "B32" "InjectiveProtocol"
Set-TextCell $ws "C32" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D32" "34.62"
Set-TextCell $ws "E32" "  +1.82%  "
Set-TextCell $ws "B33" "OKB"
Set-TextCell $ws "C33" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws "D33" "51.55"
Set-TextCell $ws "E33" "  +1.95%  "
Set-TextCell $ws "D34" "5.73"
Set-TextCell $ws "E34" "  -1.67%  "
Set-TextCell $ws "D35" "0.0845"
Set-TextCell $ws "E35" "  +2.70%  "
Set-TextCell $ws "D36" "5.29"
Set-TextCell $ws "E36" "  +0.19%  "
Set-TextCell $ws "E37" "  -0.05%  "
Set-TextCell $ws "E38" "  +0.12%  "
Set-TextCell $ws "E39" "  -2.77%  "
Set-TextCell $ws "D40" "18.01"
Set-TextCell $ws "E40" "  -0.67%  "
Set-TextCell $ws "D41" "2.57"
Set-TextCell $ws "E41" "  +2.14%  "
Set-TextCell $ws "E42" "  -1.47%  "
Set-TextCell $ws "E43" "  -1.62%  "
Set-TextCell $ws "D44" "122.26"
Set-TextCell $ws "E44" "  -3.21%  "
Set-TextCell $ws "D45" "22.10"
Set-TextCell $ws "E45" "  -7.46%  "
Set-TextCell $ws "D46" "2.075.77"
Set-TextCell $ws "E46" "  -0.19%  "
Set-TextCell $ws "D47" "3.25"
Set-TextCell $ws "E47" "  -2.47%  "
Set-TextCell $ws "E48" "  -4.35%  "
Set-TextCell $ws "E49" "  +1.27%  "
Set-TextCell $ws "D50" "0.931"
Set-TextCell $ws "E50" "  -0.03%  "
Set-TextCell $ws "E51" "  +0.33%  "
